# "add link to repo"
#
# On the "Useful Links" slide, a new line linking to the project's GitHub
# repository is added right above the existing meteor.com/tracker link.

$p = $ppt.ActivePresentation

# Find the "Useful Links" slide (title == "Useful Links") instead of a
# hard-coded index, so the script is resilient to slide ordering.
$usefulLinksSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    for ($j = 1; $j -le $candidate.Shapes.Count; $j++) {
        $candidateShape = $candidate.Shapes.Item($j)
        if ($candidateShape.Type -eq 14 -and $candidateShape.TextFrame.HasText) {
            if ($candidateShape.TextFrame.TextRange.Text.TrimEnd() -eq "Useful Links") {
                $usefulLinksSlide = $candidate
            }
        }
    }
}
if ($usefulLinksSlide -eq $null) {
    $usefulLinksSlide = $p.Slides.Item($p.Slides.Count)
}

# The links live in the body/content placeholder (the non-title shape).
$linksShape = $null
for ($j = 1; $j -le $usefulLinksSlide.Shapes.Count; $j++) {
    $candidateShape = $usefulLinksSlide.Shapes.Item($j)
    if ($candidateShape.Name -like "Content Placeholder*") {
        $linksShape = $candidateShape
    }
}
if ($linksShape -eq $null) {
    $linksShape = $usefulLinksSlide.Shapes.Item(2)
}

$linksRange = $linksShape.TextFrame.TextRange

# Locate the paragraph that holds the "https://www.meteor.com/tracker" line,
# and bail out if the repo link was already added (keeps the script
# idempotent if it were ever re-applied).
$trackerParaIndex = -1
$alreadyAdded = $false
for ($i = 1; $i -le $linksRange.Paragraphs().Count; $i++) {
    $paraText = $linksRange.Paragraphs($i).Text.TrimEnd()
    if ($paraText -eq "https://www.meteor.com/tracker") {
        $trackerParaIndex = $i
    }
    if ($paraText -eq "https://github.com/ManuelDeLeon/AllAboutReactivity") {
        $alreadyAdded = $true
    }
}

if ($trackerParaIndex -gt 0 -and -not $alreadyAdded) {
    $trackerPara = $linksRange.Paragraphs($trackerParaIndex)

    # Insert the repo link as a brand-new paragraph right above it. It
    # inherits the hyperlinked run formatting (size/hlinkClick) already in
    # place on the tracker paragraph.
    $null = $trackerPara.InsertBefore("https://github.com/ManuelDeLeon/AllAboutReactivity`r")

    # The tracker paragraph shifted down by one; re-fetch it and nudge its
    # run by re-writing its first 5 characters ("https"), which splits it
    # into two runs ("https" / "://www.meteor.com/tracker") just like the
    # original edit did.
    $trackerPara = $linksRange.Paragraphs($trackerParaIndex + 1)
    $splitPoint = $linksRange.Characters($trackerPara.Start, 5)
    $splitPoint.Text = "https"
}
